$wb = $excel.ActiveWorkbook

# --- "Sezony Ekstra" sheet: fill in matchday V/VI (columns G/H) stats for
# the 2021/2022 season block (rows 17-20) ---
$wsEkstra = $wb.Worksheets.Item("Sezony Ekstra")
$wsEkstra.Range("G17").Value = 21
$wsEkstra.Range("H17").Value = 21
$wsEkstra.Range("G18").Value = 1
$wsEkstra.Range("H18").Value = 4
$wsEkstra.Range("G19").Value = 3
$wsEkstra.Range("H19").Value = 3
$wsEkstra.Range("G20").Value = 5
$wsEkstra.Range("H20").Value = 2

# --- "Sezony I liga" sheet: fill in matchday V/VI (columns G/H) stats for
# the 2021/2022 season block (rows 17-20) ---
$wsILiga = $wb.Worksheets.Item("Sezony I liga")
$wsILiga.Range("G17").Value = 19
$wsILiga.Range("H17").Value = 26
$wsILiga.Range("G18").Value = 1
$wsILiga.Range("H18").Value = 1
$wsILiga.Range("G19").Value = 5
$wsILiga.Range("H19").Value = 5
$wsILiga.Range("G20").Value = 7
$wsILiga.Range("H20").Value = 4

# --- Update on-sheet selection / scroll position to match the author's
# last-saved cursor position on each sheet. "Sezony Ekstra" keeps focus
# only momentarily; "Sezony I liga" is the sheet left active, matching
# the workbook's activeTab. ---
$wsEkstra.Range("H21").Select()
$wsILiga.Range("H21").Select()
